$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets contain identical data tables that
# need the same four cell updates.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: cover image URL updated
    $ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg"

    # Row 4: "想去人数" (want-to-go count) incremented
    $ws.Range("F4").Value = 1503

    # Row 8: "想去人数" incremented
    $ws.Range("F8").Value = 48

    # Row 9: "想去人数" incremented
    $ws.Range("F9").Value = 307
}
